$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update styles: add yellow highlight fill for PREMIUM = Yes rows ---

# --- Column width adjustments ---
$ws.Columns.Item(3).ColumnWidth = 85.16666666666667
$ws.Columns.Item(4).ColumnWidth = 55.166666666666664
$ws.Columns.Item(8).ColumnWidth = 41.166666666666664

# --- Data rows (2-16) ---
# Row 2
$ws.Cells.Item(2, 1).Value = '1327496'
$ws.Cells.Item(2, 2).Value = 'https://aiesec.org/opportunity/global-talent/1327496'
$ws.Cells.Item(2, 3).Value = '[EXP BE] CAL - Aviation Finance Commercial Intern [EU Only]'
$ws.Cells.Item(2, 4).Value = 'Brussels, Belgium'
$ws.Cells.Item(2, 5).Value = 'Yes'
$ws.Cells.Item(2, 6).Value = '36 applicants'
$ws.Cells.Item(2, 7).Value = '6 - 18 Months'
$ws.Cells.Item(2, 8).Value = 'DHL Group'
$ws.Range("E2").Interior.Color = 65535

# Row 3
$ws.Cells.Item(3, 1).Value = '1329279'
$ws.Cells.Item(3, 2).Value = 'https://aiesec.org/opportunity/global-talent/1329279'
$ws.Cells.Item(3, 3).Value = 'Markets Commercial Ops trainee'
$ws.Cells.Item(3, 4).Value = 'Bruxelles, Belgio'
$ws.Cells.Item(3, 5).Value = 'No'
$ws.Cells.Item(3, 6).Value = '4 applicants'
$ws.Cells.Item(3, 7).Value = '6 - 18 Months'
$ws.Cells.Item(3, 8).Value = 'UCB'

# Row 4
$ws.Cells.Item(4, 1).Value = '1328731'
$ws.Cells.Item(4, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328731'
$ws.Cells.Item(4, 3).Value = 'Power Electronics Internship involving Inverter Experimentation and PSIM Simulation'
$ws.Cells.Item(4, 4).Value = '日本、兵庫県神戸市'
$ws.Cells.Item(4, 5).Value = 'No'
$ws.Cells.Item(4, 6).Value = '14 applicants'
$ws.Cells.Item(4, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(4, 8).Value = 'Sohatsu Systems Laboratory Inc.'

# Row 5
$ws.Cells.Item(5, 1).Value = '1328206'
$ws.Cells.Item(5, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328206'
$ws.Cells.Item(5, 3).Value = 'Power BI Specialist'
$ws.Cells.Item(5, 4).Value = 'Frankfurt am Main, Deutschland'
$ws.Cells.Item(5, 5).Value = 'No'
$ws.Cells.Item(5, 6).Value = '132 applicants'
$ws.Cells.Item(5, 7).Value = '6 - 18 Months'
$ws.Cells.Item(5, 8).Value = 'Greyfood GmbH'

# Row 6
$ws.Cells.Item(6, 1).Value = '1328204'
$ws.Cells.Item(6, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328204'
$ws.Cells.Item(6, 3).Value = 'Sales Support'
$ws.Cells.Item(6, 4).Value = 'Frankfurt am Main, Deutschland'
$ws.Cells.Item(6, 5).Value = 'No'
$ws.Cells.Item(6, 6).Value = '259 applicants'
$ws.Cells.Item(6, 7).Value = '3 - 6 Months'
$ws.Cells.Item(6, 8).Value = 'Greyfood GmbH'

# Row 7
$ws.Cells.Item(7, 1).Value = '1328185'
$ws.Cells.Item(7, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328185'
$ws.Cells.Item(7, 3).Value = 'Data Scientist'
$ws.Cells.Item(7, 4).Value = 'Frankfurt am Main, Deutschland'
$ws.Cells.Item(7, 5).Value = 'No'
$ws.Cells.Item(7, 6).Value = '178 applicants'
$ws.Cells.Item(7, 7).Value = '3 - 6 Months'
$ws.Cells.Item(7, 8).Value = 'Greyfood GmbH'

# Row 8
$ws.Cells.Item(8, 1).Value = '1328023'
$ws.Cells.Item(8, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328023'
$ws.Cells.Item(8, 3).Value = 'Marketing - Intern'
$ws.Cells.Item(8, 4).Value = 'Nugegoda, Sri Lanka'
$ws.Cells.Item(8, 5).Value = 'No'
$ws.Cells.Item(8, 6).Value = '27 applicants'
$ws.Cells.Item(8, 7).Value = '3 - 6 Months'
$ws.Cells.Item(8, 8).Value = 'Raffles Consolidated Pvt Ltd'

# Row 9
$ws.Cells.Item(9, 1).Value = '1328019'
$ws.Cells.Item(9, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328019'
$ws.Cells.Item(9, 3).Value = 'Kitchen Operations - Intern'
$ws.Cells.Item(9, 4).Value = 'Nugegoda, Sri Lanka'
$ws.Cells.Item(9, 5).Value = 'No'
$ws.Cells.Item(9, 6).Value = '9 applicants'
$ws.Cells.Item(9, 7).Value = '3 - 6 Months'
$ws.Cells.Item(9, 8).Value = 'Raffles Consolidated Pvt Ltd'

# Row 10
$ws.Cells.Item(10, 1).Value = '1327922'
$ws.Cells.Item(10, 2).Value = 'https://aiesec.org/opportunity/global-talent/1327922'
$ws.Cells.Item(10, 3).Value = 'Digital Marketing Intern'
$ws.Cells.Item(10, 4).Value = 'Nugegoda, Sri Lanka'
$ws.Cells.Item(10, 5).Value = 'No'
$ws.Cells.Item(10, 6).Value = '18 applicants'
$ws.Cells.Item(10, 7).Value = '3 - 6 Months'
$ws.Cells.Item(10, 8).Value = 'Starbeans Ceylon (Pvt ) Ltd'

# Row 11
$ws.Cells.Item(11, 1).Value = '1327884'
$ws.Cells.Item(11, 2).Value = 'https://aiesec.org/opportunity/global-talent/1327884'
$ws.Cells.Item(11, 3).Value = '[EXP] Customer Emission Reporting (EU Preferred)'
$ws.Cells.Item(11, 4).Value = 'Fritz-Erler-Straße 5, 53113 Bonn, Germany'
$ws.Cells.Item(11, 5).Value = 'Yes'
$ws.Cells.Item(11, 6).Value = '68 applicants'
$ws.Cells.Item(11, 7).Value = '6 - 18 Months'
$ws.Cells.Item(11, 8).Value = 'DHL Group'
$ws.Range("E11").Interior.Color = 65535

# Row 12
$ws.Cells.Item(12, 1).Value = '1327498'
$ws.Cells.Item(12, 2).Value = 'https://aiesec.org/opportunity/global-talent/1327498'
$ws.Cells.Item(12, 3).Value = '[EXP BE] ACS - Aviation Finance Commercial Intern [EU Only]'
$ws.Cells.Item(12, 4).Value = 'Brussels, Belgium'
$ws.Cells.Item(12, 5).Value = 'Yes'
$ws.Cells.Item(12, 6).Value = '59 applicants'
$ws.Cells.Item(12, 7).Value = '6 - 18 Months'
$ws.Cells.Item(12, 8).Value = 'DHL Group'
$ws.Range("E12").Interior.Color = 65535

# Row 13
$ws.Cells.Item(13, 1).Value = '1327495'
$ws.Cells.Item(13, 2).Value = 'https://aiesec.org/opportunity/global-talent/1327495'
$ws.Cells.Item(13, 3).Value = 'Content Creator'
$ws.Cells.Item(13, 4).Value = 'Sheraton Al Matar, El Nozha, Cairo Governorate, Egypt'
$ws.Cells.Item(13, 5).Value = 'No'
$ws.Cells.Item(13, 6).Value = '7 applicants'
$ws.Cells.Item(13, 7).Value = '6 - 18 Months'
$ws.Cells.Item(13, 8).Value = 'Skyline Egypt Tours'

# Row 14
$ws.Cells.Item(14, 1).Value = '1326381'
$ws.Cells.Item(14, 2).Value = 'https://aiesec.org/opportunity/global-talent/1326381'
$ws.Cells.Item(14, 3).Value = 'Business Developer'
$ws.Cells.Item(14, 4).Value = 'Sheraton Al Matar, El Nozha, Cairo Governorate, Egypt'
$ws.Cells.Item(14, 5).Value = 'No'
$ws.Cells.Item(14, 6).Value = '5 applicants'
$ws.Cells.Item(14, 7).Value = '3 - 6 Months'
$ws.Cells.Item(14, 8).Value = '5 applicants'

# Row 15
$ws.Cells.Item(15, 1).Value = '1325833'
$ws.Cells.Item(15, 2).Value = 'https://aiesec.org/opportunity/global-talent/1325833'
$ws.Cells.Item(15, 3).Value = 'Finance and Accounting Intern'
$ws.Cells.Item(15, 4).Value = 'Makati City, Metro Manila, Philippines'
$ws.Cells.Item(15, 5).Value = 'No'
$ws.Cells.Item(15, 6).Value = '21 applicants'
$ws.Cells.Item(15, 7).Value = '3 - 6 Months'
$ws.Cells.Item(15, 8).Value = 'Consistent Frozen Solutions Corporation'

# Row 16
$ws.Cells.Item(16, 1).Value = '1325830'
$ws.Cells.Item(16, 2).Value = 'https://aiesec.org/opportunity/global-talent/1325830'
$ws.Cells.Item(16, 3).Value = 'Marketing Intern'
$ws.Cells.Item(16, 4).Value = 'Makati City, Metro Manila, Philippines'
$ws.Cells.Item(16, 5).Value = 'No'
$ws.Cells.Item(16, 6).Value = '26 applicants'
$ws.Cells.Item(16, 7).Value = '3 - 6 Months'
$ws.Cells.Item(16, 8).Value = 'Consistent Frozen Solutions Corporation'

